# Weekly data refresh: insert a new observation as row 365, pushing the
# existing rows 365-393 down to 366-394 (dimension grows from T393 to T394).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 365; everything below (365-393) shifts down
# to 366-394, inheriting the formatting (incl. the date style) of the row
# it displaces.
$ws.Rows.Item(365).Insert()

# Populate the newly inserted row 365 with the new weekly record.
$ws.Range("A365").Value = 6
$ws.Range("B365").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C365").Value = "Metropolitana"
$ws.Range("D365").Value = 45265
$ws.Range("E365").Value = 13
$ws.Range("F365").Value = "Fruta"
$ws.Range("G365").Value = 100101
$ws.Range("H365").Value = "Berries"
$ws.Range("I365").Value = 100101004
$ws.Range("J365").Value = "Frambuesa"
$ws.Range("K365").Value = "Sin especificar"
$ws.Range("L365").Value = "Primera"
$ws.Range("M365").Value = 400
$ws.Range("N365").Value = 12000
$ws.Range("O365").Value = 12000
$ws.Range("P365").Value = 12000
$ws.Range("Q365").Value = "$/bandeja 2 kilos"
$ws.Range("R365").Value = "Provincia de Curicó"
$ws.Range("S365").Value = 6000
$ws.Range("T365").Value = 2
